# feat: add 2022-Q1 data
#
# Before: sheets are 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计 ("Total").
# After:  sheets are 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计.
#
#   1. The existing "总计" sheet becomes "2022-Q1" (same tab position /
#      identity) and is repopulated with the 2022-Q1 per-fund holdings
#      (same 基金代码/基金名称/.../仓位排名 layout used by the other
#      quarterly sheets).
#   2. A brand-new "总计" sheet is added right after "2022-Q1", using the
#      same 日期/持有数量(只)/持有市值(亿元) summary layout as before,
#      with a new first row for 2022-Q1 and the older rows shifted down.

$wb = $excel.ActiveWorkbook

# Applies the bold/centered/thin-box style (used for header rows and the
# column-A index markers throughout this workbook) to a single cell.
# NOTE: must be called per-cell, not on a multi-cell range at once -- a
# range-level Borders.Item(..).LineStyle only draws the OUTER edge of the
# whole range (matches real Excel semantics), not a box around every cell.
function Style-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

function Style-HeaderRange($range) {
    foreach ($cell in $range.Cells) {
        Style-HeaderCell $cell
    }
}

# ---------------------------------------------------------------------
# Step 1: rename "总计" -> "2022-Q1" and replace its contents.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Drop the old A1:D5 "总计" summary content before laying out the new
# A1:H5 per-fund table over it.
$q1.Range("A1:H5").Clear()

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
Style-HeaderRange $q1.Range("B1:H1")

# Column A index markers (styled like the header row, as in the other
# quarterly sheets).
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3
Style-HeaderRange $q1.Range("A2:A5")

# B:G hold text (fund codes keep leading zeros, decimal values keep
# trailing zeros) exactly as typed, so force text format before writing.
$q1.Range("B2:G5").NumberFormat = "@"

# Row 2: 008988 大成科技创新混合A
$q1.Range("B2").Value = "008988"
$q1.Range("C2").Value = "大成科技创新混合A"
$q1.Range("D2").Value = "2.65"
$q1.Range("E2").Value = "91.66"
$q1.Range("F2").Value = "6.49"
$q1.Range("G2").Value = "0.1720"
$q1.Range("H2").Value = 3

# Row 3: 008989 大成科技创新混合C
$q1.Range("B3").Value = "008989"
$q1.Range("C3").Value = "大成科技创新混合C"
$q1.Range("D3").Value = "1.11"
$q1.Range("E3").Value = "91.66"
$q1.Range("F3").Value = "6.49"
$q1.Range("G3").Value = "0.0720"
$q1.Range("H3").Value = 3

# Row 4: 090009 大成行业轮动混合
$q1.Range("B4").Value = "090009"
$q1.Range("C4").Value = "大成行业轮动混合"
$q1.Range("D4").Value = "1.71"
$q1.Range("E4").Value = "82.88"
$q1.Range("F4").Value = "3.16"
$q1.Range("G4").Value = "0.0540"
$q1.Range("H4").Value = 6

# Row 5: 003704 光大保德信事件驱动灵活配置混合
$q1.Range("B5").Value = "003704"
$q1.Range("C5").Value = "光大保德信事件驱动灵活配置混合"
$q1.Range("D5").Value = "3.06"
$q1.Range("E5").Value = "23.55"
$q1.Range("F5").Value = "1.31"
$q1.Range("G5").Value = "0.0401"
$q1.Range("H5").Value = 5

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the page-margin layout used by the other worksheets in this
# workbook (0.75in sides, 1in top/bottom, 0.5in header/footer).
# PageSetup margins are in points (72pt = 1in).
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"
Style-HeaderRange $total.Range("B1:D1")

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
Style-HeaderRange $total.Range("A2:A6")

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.34

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 2.16

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.89

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.07000000000000001

$q1.Range("A1").Select()
